$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.030.16'
$ws.Range("E2").Value = '  -2.96%  '
$ws.Range("D3").Value = '1.719.87'
$ws.Range("E3").Value = '  -2.78%  '
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'315.68"
$ws.Range("E5").Value = '  -3.68%  '
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").Value = "'0.4632"
$ws.Range("E7").Value = '  +3.41%  '
$ws.Range("D8").Value = "'0.3441"
$ws.Range("E8").Value = '  -3.53%  '
$ws.Range("D9").Value = "'42.34"
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").Value = "'0.07301"
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("D11").Value = "'1.052"
$ws.Range("E11").Value = '  -3.75%  '
$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = "'19.89"
$ws.Range("E13").Value = '  -4.65%  '
$ws.Range("D14").Value = "'5.881"
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("D15").Value = '1.719.65'
$ws.Range("E15").Value = '  -3.03%  '
$ws.Range("D16").Value = "'6.919"
$ws.Range("E16").Value = '  -4.07%  '
$ws.Range("D17").Value = "'89.53"
$ws.Range("E17").Value = '  -3.40%  '
$ws.Range("D18").Value = "'0.00001047"
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").Value = "'0.06324"
$ws.Range("E19").Value = '  -1.34%  '
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").Value = "'16.54"
$ws.Range("E21").Value = '  -3.71%  '
$ws.Range("D22").Value = "'5.636"
$ws.Range("E22").Value = '  -3.09%  '
$ws.Range("D23").Value = '27.084.22'
$ws.Range("E23").Value = '  -2.84%  '
$ws.Range("D24").Value = "'10.82"
$ws.Range("E24").Value = '  -4.28%  '
$ws.Range("D25").Value = "'2.134"
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").Value = "'157.02"
$ws.Range("E26").Value = '  -3.36%  '
$ws.Range("D27").Value = "'19.47"
$ws.Range("E27").Value = '  -3.75%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = "'2.156"
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = "'119.22"
$ws.Range("E29").Value = '  -4.84%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = "'1.027"
$ws.Range("E30").Value = '  -6.47%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = "'0.09093"
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").Value = "'3.597"
$ws.Range("E32").Value = '  -1.06%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = "'5.323"
$ws.Range("E33").Value = '  -4.47%  '
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").Value = "'0.02202"
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = "'11.15"
$ws.Range("E35").Value = '  -5.43%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = "'0.05840"
$ws.Range("E36").Value = '  -4.10%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = "'4.779"
$ws.Range("E37").Value = '  -3.49%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'0.1995"
$ws.Range("E38").Value = '  -4.80%  '
$ws.Range("B39").Value = 'WEMIXTOKEN'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").Value = "'1.403"
$ws.Range("E39").Value = '  +0.64%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.5976"
$ws.Range("E40").Value = '  -5.43%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'1.126"
$ws.Range("E41").Value = '  -4.84%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'7.516"
$ws.Range("E42").Value = '  -4.97%  '
$ws.Range("B43").Value = 'PancakeSwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D43").Value = "'3.625"
$ws.Range("E43").Value = '  -3.08%  '
$ws.Range("D44").Value = "'12.57"
$ws.Range("E44").Value = '  -4.75%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.5637"
$ws.Range("E45").Value = '  -3.82%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = "'119.44"
$ws.Range("E46").Value = '  -2.44%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'1.863"
$ws.Range("E47").Value = '  -4.66%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = "'0.06664"
$ws.Range("E48").Value = '  -3.53%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = "'1.088"
$ws.Range("E49").Value = '  -4.36%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = "'1.004"
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = "'69.73"
$ws.Range("E51").Value = '  -4.23%  '
